# Atualização automática: 2025-09-02 09:00:26
#
# Applies the data refresh to Sheet1:
#  - rows 7..11 rotate (old row 11 becomes the new row 7, rows 7..10 shift
#    down to become rows 8..11) — only the columns that actually differ
#    between consecutive records (A, D, E, F, G, H, I, J) are touched;
#  - row 18 gets an updated detection image / bounding box / confidence;
#  - a brand new detection record is appended as row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

function Set-NumberCell($row, $col, $number) {
    $ws.Cells.Item($row, $col).Value = $number
}

# ---- Row 7 (now holds what used to be row 11's record) ----
Set-TextCell 7 1 "2117575c-4ae1-458c-b88a-fc40f40debdb"
Set-TextCell 7 4 "image_20250727074723_ppp0.jpg"
Set-TextCell 7 5 "PLACA_20250723145134"
Set-TextCell 7 6 "Moura"
Set-NumberCell 7 7 38.06587
Set-NumberCell 7 8 -7.221796
Set-TextCell 7 9 "1490,161,1563,258"
Set-TextCell 7 10 "0.62"

# ---- Row 8 (now holds what used to be row 7's record) ----
Set-TextCell 8 1 "283b6eda-9c83-4cdd-9524-c7c394f2dc89"
Set-TextCell 8 9 "962,713,1006,765"
Set-TextCell 8 10 "0.76"

# ---- Row 9 (now holds what used to be row 8's record) ----
Set-TextCell 9 1 "a19b65d1-6f97-4841-9e1c-7446a9be92b6"
Set-TextCell 9 9 "967,614,1002,659"
Set-TextCell 9 10 "0.73"

# ---- Row 10 (now holds what used to be row 9's record) ----
Set-TextCell 10 1 "4be1b1cf-d480-453e-b5fb-d4ecd6764c4d"
Set-TextCell 10 9 "702,633,740,690"
Set-TextCell 10 10 "0.72"

# ---- Row 11 (now holds what used to be row 10's record) ----
Set-TextCell 11 1 "dfd476d4-7689-4671-a076-78fe3ce806bb"
Set-TextCell 11 9 "1254,850,1294,895"
Set-TextCell 11 10 "0.67"

# ---- Row 18: refreshed detection image / bbox / confidence ----
Set-TextCell 18 4 "image_20250808221835_ppp0.jpg"
Set-TextCell 18 9 "1182,405,1231,455"
Set-TextCell 18 10 "0.76"

# ---- Row 24: brand new detection appended at the end ----
Set-TextCell 24 1 "687a4eaa-64d4-4e21-a791-c5a0b5673343"
Set-TextCell 24 2 "mosca"

$c24 = $ws.Cells.Item(24, 3)
$c24.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$c24.Value = 45902

Set-TextCell 24 4 "image_20250902091301_ppp0.jpg"
Set-TextCell 24 5 "PLACA_20250717165933"
Set-TextCell 24 6 "Beja"
Set-NumberCell 24 7 38.02035
Set-NumberCell 24 8 -7.94715
Set-TextCell 24 9 "1,0,703,1072"
Set-TextCell 24 10 "0.66"
